$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "'71.492.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.58%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.811.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.04%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'703.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +5.86%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'174.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +4.14%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.810.78"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.05%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.03%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.163"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.96%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'7.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +5.03%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.18%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000260"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +7.34%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'36.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.61%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.458.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.818.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.54%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'71.422.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.71%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'17.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.37%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.62%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.10%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'10.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +5.72%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'483.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.75%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.25%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'84.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.91%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -2.14%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'12.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.12%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'10.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.60%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +1.75%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'3.964.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.32%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.01%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'3.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +12.18%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'7.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.19%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'2.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.84%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +6.82%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'29.64"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.49%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'9.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.38%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.11%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +1.36%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +4.33%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'6.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.68%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +10.90%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.991"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.23%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.18%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D45").Value = "'0.000317"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +15.35%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'164.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +3.64%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'45.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.09%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'48.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.93%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "'1.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.16%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +0.43%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "'417.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +6.55%  "
$ws.Range("E51").Style = "Normal"
